# Refresh cryptos list: update prices / 1h-volume %s, and fix the
# ShibaInu / WrappedEther row order (rows 16-17) to match the latest
# coinranking.com snapshot.
#
# Note: several "Price" (column D) values look like plain numbers once
# updated (e.g. "9.79", "1.00", "13.90") but must stay TEXT so trailing
# zeros / the dotted-thousands formatting of the source data survive.
# A leading apostrophe (the classic Excel "force text" quote-prefix)
# is used for those so Excel doesn't silently coerce them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.150.15'
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").Value = '3.522.94'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''584.06'
$ws.Range("E5").Value = '  -1.20%  '
$ws.Range("D6").Value = '''134.25'
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("D7").Value = '3.522.28'
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("D11").Value = '''7.13'
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("D12").Value = '''0.377'
$ws.Range("E12").Value = '  -2.06%  '
$ws.Range("D13").Value = '4.124.44'
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").Value = '''27.45'
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("D15").Value = '''0.119'
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.531.56'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.0000179'
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("D18").Value = '64.203.09'
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("D19").Value = '''9.79'
$ws.Range("E19").Value = '  -2.52%  '
$ws.Range("D20").Value = '''13.90'
$ws.Range("E20").Value = '  -2.89%  '
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("D22").Value = '''382.41'
$ws.Range("E22").Value = '  -2.17%  '
$ws.Range("D23").Value = '''0.572'
$ws.Range("E23").Value = '  -0.89%  '
$ws.Range("D24").Value = '3.665.68'
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").Value = '''74.07'
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '''5.61'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("E28").Value = '  +3.56%  '
$ws.Range("D29").Value = '''1.58'
$ws.Range("E29").Value = '  -2.21%  '
$ws.Range("D30").Value = '''7.45'
$ws.Range("E30").Value = '  -2.39%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").Value = '''8.41'
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("D34").Value = '3.537.00'
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D36").Value = '''23.58'
$ws.Range("E36").Value = '  -1.97%  '
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("D38").Value = '''5.40'
$ws.Range("E38").Value = '  +3.19%  '
$ws.Range("E39").Value = '  +0.53%  '
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("D41").Value = '''160.77'
$ws.Range("E41").Value = '  -5.14%  '
$ws.Range("D42").Value = '''0.0785'
$ws.Range("E42").Value = '  -2.34%  '
$ws.Range("D43").Value = '''26.62'
$ws.Range("E43").Value = '  +2.29%  '
$ws.Range("D44").Value = '''0.812'
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = '''1.21'
$ws.Range("E46").Value = '  -2.98%  '
$ws.Range("D47").Value = '''41.66'
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("D50").Value = '2.483.02'
$ws.Range("D51").Value = '''6.80'
$ws.Range("E51").Value = '  -1.08%  '
